$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H10").Value = 1200
$ws.Range("J10").Value = 1200
$ws.Range("L10").Value = 1200
$ws.Range("N10").Value = -1786
$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").ClearContents()
$ws.Range("H137").Value = 23600
$ws.Range("I137").Value = 31500
$ws.Range("K137").Value = 94500
$ws.Range("M137").Value = -91950
$ws.Range("H138").Value = 7336.4404
$ws.Range("I138").Value = 6778.2593
$ws.Range("J138").Value = 7600.8423
$ws.Range("K138").Value = 20334.7779
$ws.Range("L138").Value = 22802.5269
$ws.Range("M138").Value = -15194.7779
$ws.Range("N138").Value = -33082.5269

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3127.75
$ws.Range("I45").Value = 2383.5
$ws.Range("K45").Value = 2383.5
$ws.Range("M45").Value = -2006.5
$ws.Range("H61").Value = 2088.8
$ws.Range("I61").Value = 2098.6667
$ws.Range("J61").Value = 2000
$ws.Range("K61").Value = 2098.6667
$ws.Range("L61").Value = 2000
$ws.Range("M61").Value = -1886.6667
$ws.Range("N61").Value = -2424
$ws.Range("H74").Value = 7332.8
$ws.Range("I74").Value = 2000
$ws.Range("K74").Value = 2000
$ws.Range("M74").Value = -1126
$ws.Range("H77").Value = 7332.8
$ws.Range("I77").Value = 2000
$ws.Range("K77").Value = 10000
$ws.Range("M77").Value = -5632
$ws.Range("H132").Value = 4736.923
$ws.Range("I132").Value = 2393.111
$ws.Range("J132").Value = 10010.5
$ws.Range("K132").Value = 7179.333
$ws.Range("L132").Value = 30031.5
$ws.Range("M132").Value = -4649.333
$ws.Range("N132").Value = -35091.5
$ws.Range("H136").Value = 2088.8
$ws.Range("I136").Value = 2098.6667
$ws.Range("J136").Value = 2000
$ws.Range("K136").Value = 6296.000100000001
$ws.Range("L136").Value = 6000
$ws.Range("M136").Value = -3746.000100000001
$ws.Range("N136").Value = -11100

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H75").Value = 33249
$ws.Range("I75").Value = 6499
$ws.Range("J75").Value = 59999
$ws.Range("K75").Value = 6499
$ws.Range("L75").Value = 59999
$ws.Range("M75").Value = -5563
$ws.Range("N75").Value = -61871
$ws.Range("H78").Value = 33249
$ws.Range("I78").Value = 6499
$ws.Range("J78").Value = 59999
$ws.Range("K78").Value = 19497
$ws.Range("L78").Value = 179997
$ws.Range("M78").Value = -14817
$ws.Range("N78").Value = -189357
$ws.Range("H134").Value = 3220.6924
$ws.Range("I134").Value = 1426.375
$ws.Range("J134").Value = 6091.6
$ws.Range("K134").Value = 4279.125
$ws.Range("L134").Value = 18274.8
$ws.Range("M134").Value = -1744.125
$ws.Range("N134").Value = -23344.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H13").Value = 2500
$ws.Range("I13").Value = 2500
$ws.Range("K13").Value = 2500
$ws.Range("M13").Value = -2361
$ws.Range("H31").Value = 4244.4375
$ws.Range("I31").Value = 2613.7778
$ws.Range("J31").Value = 6341
$ws.Range("K31").Value = 2613.7778
$ws.Range("L31").Value = 6341
$ws.Range("M31").Value = -2318.7778
$ws.Range("N31").Value = -6931
$ws.Range("H34").Value = 4244.4375
$ws.Range("I34").Value = 2613.7778
$ws.Range("J34").Value = 6341
$ws.Range("K34").Value = 2613.7778
$ws.Range("L34").Value = 6341
$ws.Range("M34").Value = -2411.7778
$ws.Range("N34").Value = -6745
$ws.Range("H58").Value = 5142.9
$ws.Range("I58").Value = 1398.6
$ws.Range("J58").Value = 8887.200000000001
$ws.Range("K58").Value = 1398.6
$ws.Range("L58").Value = 8887.200000000001
$ws.Range("M58").Value = -1195.6
$ws.Range("N58").Value = -9293.200000000001
$ws.Range("H99").Value = 14888.723
$ws.Range("I99").Value = 12834.167
$ws.Range("J99").Value = 15916
$ws.Range("K99").Value = 12834.167
$ws.Range("L99").Value = 15916
$ws.Range("M99").Value = -11336.167
$ws.Range("N99").Value = -18912
$ws.Range("H126").Value = 14888.723
$ws.Range("I126").Value = 12834.167
$ws.Range("J126").Value = 15916
$ws.Range("K126").Value = 38502.501
$ws.Range("L126").Value = 47748
$ws.Range("M126").Value = -36032.501
$ws.Range("N126").Value = -52688
$ws.Range("H132").Value = 2988.1365
$ws.Range("I132").Value = 2968.125
$ws.Range("J132").Value = 3041.5
$ws.Range("K132").Value = 8904.375
$ws.Range("L132").Value = 9124.5
$ws.Range("M132").Value = -6374.375
$ws.Range("N132").Value = -14184.5
$ws.Range("H136").Value = 5142.9
$ws.Range("I136").Value = 1398.6
$ws.Range("J136").Value = 8887.200000000001
$ws.Range("K136").Value = 4195.799999999999
$ws.Range("L136").Value = 26661.6
$ws.Range("M136").Value = -1645.799999999999
$ws.Range("N136").Value = -31761.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 1741.6364
$ws.Range("I14").Value = 1741.6364
$ws.Range("K14").Value = 5224.9092
$ws.Range("M14").Value = -5051.9092
$ws.Range("H34").Value = 2466.3333
$ws.Range("J34").Value = 5000
$ws.Range("L34").Value = 15000
$ws.Range("N34").Value = -15168
$ws.Range("H112").Value = 3604.5
$ws.Range("I112").Value = 3342.3333
$ws.Range("J112").Value = 3866.6667
$ws.Range("K112").Value = 10026.9999
$ws.Range("L112").Value = 11600.0001
$ws.Range("M112").Value = -8918.999899999999
$ws.Range("N112").Value = -13816.0001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3806
$ws.Range("I132").Value = 2649.3
$ws.Range("J132").Value = 8432.799999999999
$ws.Range("K132").Value = 7947.900000000001
$ws.Range("L132").Value = 25298.4
$ws.Range("M132").Value = -5417.900000000001
$ws.Range("N132").Value = -30358.4
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H74").Value = 50000
$ws.Range("I74").Value = 50000
$ws.Range("K74").Value = 50000
$ws.Range("M74").Value = -49002
$ws.Range("H77").Value = 50000
$ws.Range("I77").Value = 50000
$ws.Range("K77").Value = 150000
$ws.Range("M77").Value = -145008
$ws.Range("H82").Value = 2077.125
$ws.Range("I82").Value = 2736.125
$ws.Range("J82").Value = 1418.125
$ws.Range("K82").Value = 2736.125
$ws.Range("L82").Value = 1418.125
$ws.Range("M82").Value = -2375.125
$ws.Range("N82").Value = -2140.125
$ws.Range("H85").Value = 2077.125
$ws.Range("I85").Value = 2736.125
$ws.Range("J85").Value = 1418.125
$ws.Range("K85").Value = 2736.125
$ws.Range("L85").Value = 1418.125
$ws.Range("M85").Value = -1488.125
$ws.Range("N85").Value = -3914.125
$ws.Range("H132").Value = 3552.158
$ws.Range("I132").Value = 1888.2727
$ws.Range("J132").Value = 5840
$ws.Range("K132").Value = 5664.8181
$ws.Range("L132").Value = 17520
$ws.Range("M132").Value = -3134.8181
$ws.Range("N132").Value = -22580

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 6591.5
$ws.Range("J62").Value = 6870.1
$ws.Range("L62").Value = 6870.1
$ws.Range("N62").Value = -8118.1
$ws.Range("H65").Value = 6591.5
$ws.Range("J65").Value = 6870.1
$ws.Range("L65").Value = 34350.5
$ws.Range("N65").Value = -40590.5
$ws.Range("H113").Value = 996.2759
$ws.Range("I113").Value = 633.6667
$ws.Range("J113").Value = 1252.2354
$ws.Range("K113").Value = 1901.0001
$ws.Range("L113").Value = 3756.7062
$ws.Range("M113").Value = 268.9999
$ws.Range("N113").Value = -8096.706200000001
$ws.Range("H132").Value = 2230.375
$ws.Range("I132").Value = 1187.8889
$ws.Range("J132").Value = 3570.7144
$ws.Range("K132").Value = 3563.6667
$ws.Range("L132").Value = 10712.1432
$ws.Range("M132").Value = -1033.6667
$ws.Range("N132").Value = -15772.1432
$ws.Range("H136").Value = 49296.953
$ws.Range("I136").Value = 926.2941
$ws.Range("J136").Value = 254872.25
$ws.Range("K136").Value = 2778.8823
$ws.Range("L136").Value = 764616.75
$ws.Range("M136").Value = -228.8822999999998
$ws.Range("N136").Value = -769716.75
